$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 used to hold a formula "(B2 + C2)"; it's now a plain text label "x_sum"
$ws.Range("D1").Value = "x_sum"

# D2:D4 now all sum A2:A4
$ws.Range("D2").Formula = "=SUM(A2:A4)"
$ws.Range("D3").Formula = "=SUM(A2:A4)"
$ws.Range("D4").Formula = "=SUM(A2:A4)"

# D5:D7 now all sum A5:A7 (D7 is a brand new cell)
$ws.Range("D5").Formula = "=SUM(A5:A7)"
$ws.Range("D6").Formula = "=SUM(A5:A7)"
$ws.Range("D7").Formula = "=SUM(A5:A7)"
